$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new blank rows before the existing data (existing rows 2-8 shift to 12-18)
$ws.Rows("2:11").Insert()

# Copy formatting (styles) from the now-shifted first data row (row 12) down onto
# the newly inserted rows so they match the rest of the dataset's look.
$ws.Range("A12:G12").Copy()
$ws.Range("A2:G11").PasteSpecial(-4122)
$ws.Rows("2:11").RowHeight = 19.5

$ws.Range("A2").Value = 25569.313026319443
$ws.Range("B2").Value = "GDHG"
$ws.Range("C2").Value = 0.1895
$ws.Range("D2").Value = 0.178
$ws.Range("E2").Value = 500
$ws.Range("F2").Value = -6.890000000000001
$ws.Range("G2").Value = -6.890000000000001
$ws.Range("A3").Value = 25569.313026319443
$ws.Range("B3").Value = "LPTV"
$ws.Range("C3").Value = 0.1506
$ws.Range("D3").Value = 0.14
$ws.Range("E3").Value = 500
$ws.Range("F3").Value = -6.440000000000012
$ws.Range("G3").Value = -13.33000000000001
$ws.Range("A4").Value = 25569.313026319443
$ws.Range("B4").Value = "SNTI"
$ws.Range("C4").Value = 0.4851
$ws.Range("D4").Value = 0.52
$ws.Range("E4").Value = 500
$ws.Range("F4").Value = 16.31000000000003
$ws.Range("G4").Value = 2.980000000000018
$ws.Range("A5").Value = 25569.313026319443
$ws.Range("B5").Value = "SNTI"
$ws.Range("C5").Value = 0.47
$ws.Range("D5").Value = 0.52
$ws.Range("E5").Value = 400
$ws.Range("F5").Value = 18.85999999999999
$ws.Range("G5").Value = 21.84
$ws.Range("A6").Value = 25569.313026319443
$ws.Range("B6").Value = "LUCY"
$ws.Range("C6").Value = 0.5165
$ws.Range("D6").Value = 0.4921
$ws.Range("E6").Value = 500
$ws.Range("F6").Value = -13.34000000000003
$ws.Range("G6").Value = 8.499999999999972
$ws.Range("A7").Value = 25569.313026319443
$ws.Range("B7").Value = "GOVX"
$ws.Range("C7").Value = 3.61
$ws.Range("D7").Value = 3.55
$ws.Range("E7").Value = 65
$ws.Range("F7").Value = -5.04000000000002
$ws.Range("G7").Value = 3.459999999999951
$ws.Range("A8").Value = 25569.31302633102
$ws.Range("B8").Value = "TSLA"
$ws.Range("C8").Value = 210.82
$ws.Range("D8").Value = 226.46
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 30.14000000000004
$ws.Range("G8").Value = 30.14000000000004
$ws.Range("A9").Value = 25569.31302633102
$ws.Range("B9").Value = "OPTT"
$ws.Range("C9").Value = 0.455
$ws.Range("D9").Value = 0.54
$ws.Range("E9").Value = 100
$ws.Range("F9").Value = 7.359999999999999
$ws.Range("G9").Value = 37.50000000000004
$ws.Range("A10").Value = 25569.31302633102
$ws.Range("B10").Value = "MBIO"
$ws.Range("C10").Value = 0.535
$ws.Range("D10").Value = 0.566
$ws.Range("E10").Value = 100
$ws.Range("F10").Value = 1.959999999999994
$ws.Range("G10").Value = 39.46000000000004
$ws.Range("A11").Value = 25569.31302633102
$ws.Range("B11").Value = "ANVS"
$ws.Range("C11").Value = 6.47
$ws.Range("D11").Value = 10.71
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 3.100000000000001
$ws.Range("G11").Value = 42.56000000000004